$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-08-18 14:49:50"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-08-18 14:49:44"
$zhcn.Range("K4").Value = "2016-08-18 14:50:17"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-08-18 14:49:50"
$dede.Range("K4").Value = "2016-08-18 14:50:30"
